$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 7999.5
$ws.Range("I13").Value = 7999
$ws.Range("K13").Value = 7999
$ws.Range("M13").Value = -7830

$ws.Range("H15").Value = 1829.6666
$ws.Range("I15").Value = 1829.6666
$ws.Range("K15").Value = 5488.9998
$ws.Range("M15").Value = -5319.9998

$ws.Range("H28").Value = 1576.7858
$ws.Range("I28").Value = 1019.6
$ws.Range("K28").Value = 1019.6
$ws.Range("M28").Value = -534.6

$ws.Range("H45").Value = 8887
$ws.Range("J45").Value = 8887
$ws.Range("L45").Value = 26661
$ws.Range("N45").Value = -27045

$ws.Range("H74").Value = 4989.6665
$ws.Range("I74").Value = 4989.6665
$ws.Range("K74").Value = 4989.6665
$ws.Range("M74").Value = -4053.6665

$ws.Range("H77").Value = 4989.6665
$ws.Range("I77").Value = 4989.6665
$ws.Range("K77").Value = 24948.3325
$ws.Range("M77").Value = -20268.3325

$ws.Range("H98").Value = 1256.5714
$ws.Range("I98").Value = 1256.5714
$ws.Range("K98").Value = 1256.5714
$ws.Range("M98").Value = 241.4286

$ws.Range("H122").Value = 1256.5714
$ws.Range("I122").Value = 1256.5714
$ws.Range("K122").Value = 3769.7142
$ws.Range("M122").Value = -1319.7142

$ws.Range("H127").Value = 944.25
$ws.Range("I127").Value = 944.25
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2832.75
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = 2127.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3727.4119
$ws.Range("I32").Value = 2672.0322
$ws.Range("K32").Value = 2672.0322
$ws.Range("M32").Value = -2385.0322

$ws.Range("H47").Value = 40000
$ws.Range("I47").Value = 40000
$ws.Range("K47").Value = 40000
$ws.Range("M47").Value = -39275

$ws.Range("H61").Value = 2111.75
$ws.Range("I61").Value = 2111.75
$ws.Range("K61").Value = 2111.75
$ws.Range("M61").Value = -1899.75

$ws.Range("H74").Value = 1294.125
$ws.Range("I74").Value = 988.5
$ws.Range("K74").Value = 988.5
$ws.Range("M74").Value = -114.5

$ws.Range("H77").Value = 1294.125
$ws.Range("I77").Value = 988.5
$ws.Range("K77").Value = 4942.5
$ws.Range("M77").Value = -574.5

$ws.Range("H110").Value = 10603.8
$ws.Range("I110").Value = 10603.8
$ws.Range("K110").Value = 10603.8
$ws.Range("M110").Value = -8558.799999999999

$ws.Range("H122").Value = 1558.5454
$ws.Range("I122").Value = 1575.619
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4726.857
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2276.857
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 3124.5
$ws.Range("I132").Value = 4033
$ws.Range("J132").Value = 399
$ws.Range("K132").Value = 12099
$ws.Range("L132").Value = 1197
$ws.Range("M132").Value = -9569
$ws.Range("N132").Value = -6257

$ws.Range("H136").Value = 2111.75
$ws.Range("I136").Value = 2111.75
$ws.Range("K136").Value = 6335.25
$ws.Range("M136").Value = -3785.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10336

$ws.Range("H20").Value = 2717
$ws.Range("J20").Value = 2767.25
$ws.Range("L20").Value = 2767.25
$ws.Range("N20").Value = -3261.25

$ws.Range("H44").Value = 49000
$ws.Range("J44").Value = 49000
$ws.Range("L44").Value = 49000
$ws.Range("N44").Value = -49994

$ws.Range("H107").Value = 1050.4546
$ws.Range("I107").Value = 907
$ws.Range("J107").Value = 1433
$ws.Range("K107").Value = 907
$ws.Range("L107").Value = 1433
$ws.Range("M107").Value = 1013
$ws.Range("N107").Value = -5273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1498.5
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 1498.5
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H58").Value = 8966.666999999999
$ws.Range("J58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10406

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

$ws.Range("H122").Value = 3244.6428
$ws.Range("I122").Value = 3499.875
$ws.Range("J122").Value = 2904.3333
$ws.Range("K122").Value = 10499.625
$ws.Range("L122").Value = 8712.999899999999
$ws.Range("M122").Value = -8049.625
$ws.Range("N122").Value = -13612.9999

$ws.Range("H132").Value = 2165
$ws.Range("I132").Value = 1997.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5992.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3462.5
$ws.Range("N132").Value = -12560

$ws.Range("H136").Value = 8966.666999999999
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 246.4
$ws.Range("I92").Value = 220.75
$ws.Range("J92").Value = 349
$ws.Range("K92").Value = 662.25
$ws.Range("L92").Value = 1047
$ws.Range("M92").Value = 585.75
$ws.Range("N92").Value = -3543

$ws.Range("H129").Value = 2554
$ws.Range("J129").Value = 2822.1667
$ws.Range("L129").Value = 8466.500100000001
$ws.Range("N129").Value = -18466.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 21999.166
$ws.Range("I44").Value = 18400
$ws.Range("K44").Value = 18400
$ws.Range("M44").Value = -17804

$ws.Range("H55").Value = 3619.8
$ws.Range("J55").Value = 3900
$ws.Range("L55").Value = 3900
$ws.Range("N55").Value = -4554

$ws.Range("H102").Value = 2337.1667
$ws.Range("I102").Value = 2341.3333
$ws.Range("J102").Value = 2333
$ws.Range("K102").Value = 2341.3333
$ws.Range("L102").Value = 2333
$ws.Range("M102").Value = -719.3332999999998
$ws.Range("N102").Value = -5577

$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws.Range("H126").Value = 1938
$ws.Range("I126").Value = 1938
$ws.Range("K126").Value = 5814
$ws.Range("M126").Value = -3344

$ws.Range("H132").Value = 1256.6666
$ws.Range("I132").Value = 1256.6666
$ws.Range("K132").Value = 3769.9998
$ws.Range("M132").Value = -1239.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5163.1113
$ws.Range("I7").Value = 5358
$ws.Range("K7").Value = 5358
$ws.Range("M7").Value = -5246

$ws.Range("H22").Value = 1186
$ws.Range("I22").Value = 1234.375
$ws.Range("J22").Value = 799
$ws.Range("K22").Value = 1234.375
$ws.Range("L22").Value = 799
$ws.Range("M22").Value = -939.375
$ws.Range("N22").Value = -1389

$ws.Range("H27").Value = 1186
$ws.Range("I27").Value = 1234.375
$ws.Range("J27").Value = 799
$ws.Range("K27").Value = 1234.375
$ws.Range("L27").Value = 799
$ws.Range("M27").Value = -1127.375
$ws.Range("N27").Value = -1013

$ws.Range("H30").Value = 930.75
$ws.Range("I30").Value = 1124.3334
$ws.Range("K30").Value = 1124.3334
$ws.Range("M30").Value = -1016.3334

$ws.Range("H126").Value = 5163.1113
$ws.Range("I126").Value = 5358
$ws.Range("K126").Value = 16074
$ws.Range("M126").Value = -13604

$ws.Range("H132").Value = 3038.2693
$ws.Range("I132").Value = 3591.6155
$ws.Range("J132").Value = 2484.923
$ws.Range("K132").Value = 10774.8465
$ws.Range("L132").Value = 7454.768999999999
$ws.Range("M132").Value = -8244.8465
$ws.Range("N132").Value = -12514.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 11747.5
$ws.Range("I8").Value = 11500
$ws.Range("J8").Value = 11995
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 11995
$ws.Range("M8").Value = -11360
$ws.Range("N8").Value = -12275

$ws.Range("H51").Value = 233331.75
$ws.Range("I51").Value = 233331.75
$ws.Range("K51").Value = 233331.75
$ws.Range("M51").Value = -232821.75

$ws.Range("H100").Value = 1599.6666
$ws.Range("I100").Value = 300
$ws.Range("K100").Value = 600
$ws.Range("M100").Value = -59

$ws.Range("H122").Value = 2724.7334
$ws.Range("I122").Value = 2743.6667
$ws.Range("K122").Value = 8231.000100000001
$ws.Range("M122").Value = -5781.000100000001

$ws.Range("H126").Value = 2243.3044
$ws.Range("I126").Value = 1995.2
$ws.Range("K126").Value = 5985.6
$ws.Range("M126").Value = -3515.6
